# Apply 2024-11-16 data refresh to violent-crime-full-year workbook.
# For each affected worksheet (citywide totals, by-neighborhood summary, and
# individual neighborhood sheets), update the 2024 (column K) crime counts, and
# for a few sheets the preceding 2023 (column J) / 2020 (column G) counts as well,
# to reflect newly-ingested incident records for that day.

$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7095
$ws.Range("K3").Value = 7354
$ws.Range("G4").Value = 1492
$ws.Range("J4").Value = 1842
$ws.Range("K4").Value = 1533
$ws.Range("K5").Value = 522
$ws.Range("K6").Value = 8099
$ws.Range("G7").Value = 24720
$ws.Range("J7").Value = 29309
$ws.Range("K7").Value = 24603

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K4").Value = 89
$ws.Range("K5").Value = 67
$ws.Range("K7").Value = 744
$ws.Range("K8").Value = 1603
$ws.Range("K10").Value = 141
$ws.Range("K19").Value = 718
$ws.Range("K20").Value = 601
$ws.Range("K21").Value = 83
$ws.Range("K23").Value = 249
$ws.Range("K24").Value = 78
$ws.Range("K27").Value = 231
$ws.Range("K29").Value = 1355
$ws.Range("K33").Value = 1050
$ws.Range("K34").Value = 138
$ws.Range("K36").Value = 314
$ws.Range("K37").Value = 831
$ws.Range("K42").Value = 908
$ws.Range("K43").Value = 201
$ws.Range("K48").Value = 314
$ws.Range("K49").Value = 137
$ws.Range("K51").Value = 310
$ws.Range("K52").Value = 640
$ws.Range("K54").Value = 479
$ws.Range("K61").Value = 20
$ws.Range("G63").Value = 295
$ws.Range("J63").Value = 122
$ws.Range("K65").Value = 574
$ws.Range("K67").Value = 965
$ws.Range("K76").Value = 332
$ws.Range("K78").Value = 296
$ws.Range("K79").Value = 606
$ws.Range("K80").Value = 91
$ws.Range("K81").Value = 18
$ws.Range("K82").Value = 30
$ws.Range("K83").Value = 523
$ws.Range("K84").Value = 194
$ws.Range("K85").Value = 1123
$ws.Range("K89").Value = 368
$ws.Range("K92").Value = 90
$ws.Range("K95").Value = 401
$ws.Range("K96").Value = 266
$ws.Range("K99").Value = 422
$ws.Range("G101").Value = 24720
$ws.Range("J101").Value = 29309
$ws.Range("K101").Value = 24603

# West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 112
$ws.Range("K7").Value = 266

# Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K6").Value = 203
$ws.Range("K7").Value = 744

# Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 153
$ws.Range("K6").Value = 155

# Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 103
$ws.Range("K3").Value = 116
$ws.Range("K4").Value = 40
$ws.Range("K7").Value = 368

# South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 368
$ws.Range("K6").Value = 278
$ws.Range("K7").Value = 1123

# Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 173
$ws.Range("K6").Value = 236
$ws.Range("K7").Value = 640

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 446
$ws.Range("K3").Value = 485
$ws.Range("K6").Value = 536
$ws.Range("K7").Value = 1603

# South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K6").Value = 122
$ws.Range("K7").Value = 523

# Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 262
$ws.Range("K7").Value = 1050

# West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 137
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 401

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 241
$ws.Range("K3").Value = 272
$ws.Range("K7").Value = 831

# New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 139
$ws.Range("K7").Value = 574

# Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 177
$ws.Range("K7").Value = 422

# North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 268
$ws.Range("K3").Value = 348
$ws.Range("K4").Value = 55
$ws.Range("K7").Value = 965

# South Deering
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 81
$ws.Range("K7").Value = 194

# Lincoln Park
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 137

# Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K4").Value = 29
$ws.Range("K7").Value = 479

# Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 484
$ws.Range("K4").Value = 63
$ws.Range("K6").Value = 394
$ws.Range("K7").Value = 1355

# Lake View
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 76
$ws.Range("K7").Value = 314

# Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 210
$ws.Range("K4").Value = 33
$ws.Range("K7").Value = 718

# River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 74
$ws.Range("K4").Value = 25
$ws.Range("K7").Value = 332

# Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 245
$ws.Range("K3").Value = 270
$ws.Range("K4").Value = 41
$ws.Range("K6").Value = 337
$ws.Range("K7").Value = 908

# Avondale
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K3").Value = 26
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 141

# Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K4").Value = 27
$ws.Range("K7").Value = 296

# Dunning
$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 78

# Douglas
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 249

# Chinatown
$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 83

# Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 201
$ws.Range("K6").Value = 155
$ws.Range("K7").Value = 606

# Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 194
$ws.Range("K7").Value = 601

# Grand Boulevard
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 119
$ws.Range("K7").Value = 314

# Garfield Ridge
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 138

# West Elsdon
$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 90

# Armour Square
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 67

# Edgewater
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 62
$ws.Range("K7").Value = 231

# Little Italy, UIC
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 86
$ws.Range("K7").Value = 310

# Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 201

# Sheffield & DePaul
$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K3").Value = 5
$ws.Range("K6").Value = 30

# Rush & Division
$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 91

# Archer Heights
$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 89

# Mount Greenwood
$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 20

# Sauganash,Forest Glen
$ws = $wb.Worksheets.Item("Sauganash,Forest Glen")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 18
